# Update a handful of computed values in row 2 of Sheet1 with refreshed
# results pulled from the server.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 35.67999920255268
$ws.Range("E2").Value = 37.01077875011297
$ws.Range("J2").Value = 37.12425180097849
$ws.Range("O2").Value = 31.0610156868254
$ws.Range("Q2").Value = 33.90553370188188
